$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row (row 7): add the required boolean element as literal text "true".
# A direct $cell.Value = "true" assignment gets auto-coerced to a native Boolean by
# Excel's type inference, which is not what the source data needs (it must stay a
# plain string, matching the rest of this metadata table). Round-tripping through a
# formula + copy/paste-values sidesteps that coercion and keeps formatting intact.
$cell = $ws.Range("B7")
$cell.Formula = '="true"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues

# "Date" row (row 8): bump the generation timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
